$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.081.14'
$ws.Range("E2").Value = '  -2.73%  '
$ws.Range("D3").Value = '2.384.08'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.50'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.82'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.536'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").Value = '2.385.01'
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("E10").Value = '  -4.44%  '
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.336'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.57'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.08%  '
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("E16").Value = '  -3.74%  '
$ws.Range("D17").Value = '60.368.16'
$ws.Range("E17").Value = '  -2.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.40'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +15.41%  '
$ws.Range("D19").Value = '2.385.78'
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.36'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.99'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.77'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.36'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '549.33'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.92'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -13.48%  '
$ws.Range("D29").Value = '2.504.20'
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("D30").Value = '0.0₃0898'
$ws.Range("E30").Value = '  -2.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.28'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.57%  '
$ws.Range("E33").Value = '  -4.24%  '
$ws.Range("E34").Value = '  -2.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '152.81'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.40'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.366'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.48'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.65%  '
$ws.Range("E40").Value = '  -0.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.97'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.27%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.10'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("E44").Value = '  -3.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.27'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.01%  '
$ws.Range("E46").Value = '  -3.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.56'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.48'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.586'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.09%  '
$ws.Range("E50").Value = '  -2.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.76'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.53%  '
